$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.437.72"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "2.276.65"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.01"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.25"
$ws.Range("E6").Value = "  +4.97%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.91"
$ws.Range("E10").Value = "  +10.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.68"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "2.634.63"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.41"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "2.285.22"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").Value = "42.331.97"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.52"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.96"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.73"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.55"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.59"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.86"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.73"
$ws.Range("E28").Value = "  +6.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.50"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.12"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.25"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0740"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.04"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.06"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.44"
$ws.Range("E42").Value = "  +14.29%  "
$ws.Range("D43").Value = "1.998.64"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  -2.77%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.97"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.85"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.14"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.92"
$ws.Range("E51").Value = "  +0.79%  "
